$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1347866666666667
$ws.Range("H2").Value = 0.40436
$ws.Range("I2").Value = 0.03419045085634245
$ws.Range("J2").Value = 0.03419045085634244
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.03911365294222222
$ws.Range("R2").Value = 0.35202287648
$ws.Range("S2").Value = 0.001172523758381379
$ws.Range("T2").Value = 0.001172523758381379

# Row 3
$ws.Range("G3").Value = 0.1347866666666667
$ws.Range("H3").Value = 0.40436
$ws.Range("I3").Value = 0.03419045085634245
$ws.Range("J3").Value = 0.03419045085634244
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 0.9583804651911111
$ws.Range("R3").Value = 8.62542418672
$ws.Range("S3").Value = 0.02872970895009768
$ws.Range("T3").Value = 0.02872970895009768

# Row 4
$ws.Range("G4").Value = 0.1347866666666667
$ws.Range("H4").Value = 0.40436
$ws.Range("I4").Value = 0.03419045085634245
$ws.Range("J4").Value = 0.03419045085634244
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 0.1430485951155556
$ws.Range("R4").Value = 1.28743735604
$ws.Range("S4").Value = 0.004288218147863386
$ws.Range("T4").Value = 0.004288218147863385

# Row 5
$ws.Range("I5").Value = 0.3318597741685039
$ws.Range("J5").Value = 0.3318597741685039
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 0.3796454187413333
$ws.Range("R5").Value = 3.416808768672
$ws.Range("S5").Value = 0.01138076450932405
$ws.Range("T5").Value = 0.01138076450932404

# Row 6
$ws.Range("I6").Value = 0.3318597741685039
$ws.Range("J6").Value = 0.3318597741685039
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.2788566539869898
$ws.Range("T6").Value = 0.2788566539869897

# Row 7
$ws.Range("I7").Value = 0.3318597741685039
$ws.Range("J7").Value = 0.3318597741685039
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("S7").Value = 0.04162235567219015
$ws.Range("T7").Value = 0.04162235567219013

# Row 8
$ws.Range("I8").Value = 0.6339497749751537
$ws.Range("J8").Value = 0.6339497749751537
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 0.7252344107822223
$ws.Range("R8").Value = 6.52710969704
$ws.Range("S8").Value = 0.02174060751354522
$ws.Range("T8").Value = 0.02174060751354522

# Row 9
$ws.Range("I9").Value = 0.6339497749751537
$ws.Range("J9").Value = 0.6339497749751537
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.5326982261960279
$ws.Range("T9").Value = 0.5326982261960279

# Row 10
$ws.Range("I10").Value = 0.6339497749751537
$ws.Range("J10").Value = 0.6339497749751537
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 2.652367032685556
$ws.Range("S10").Value = 0.0795109412655806
$ws.Range("T10").Value = 0.07951094126558057
